{"js": "// The \"COMPETENCES TECHNIQUES\" skill list got reordered. The five skill\n// lines keep their exact text and paragraph formatting; only the order in\n// which they appear changes:\n//   old: Bases de donn\u00e9es, Autres, Visualisation, ML/AI, MLOps\n//   new: Visualisation, MLOps, Autres, ML/AI, Bases de donn\u00e9es\n// Since every one of those paragraphs shares identical paragraph\n// properties, the reorder can be applied by rewriting each paragraph's\n// text in place (cheapest, most robust transform \u2014 no node moves needed).\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst oldOrder = [\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Autres : utilisables, 09\",\n  \"Visualisation : excel, tableau\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n];\nconst newOrder = [\n  \"Visualisation : excel, tableau\",\n  \"MLOps : vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Autres : utilisables, 09\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n];\n\n// Locate the contiguous run of these five paragraphs (they appear back to\n// back under \"COMPETENCES TECHNIQUES\").\nlet startIndex = -1;\nfor (let i = 0; i + oldOrder.length <= paras.items.length; i++) {\n  let matches = true;\n  for (let j = 0; j < oldOrder.length; j++) {\n    if (paras.items[i + j].text !== oldOrder[j]) {\n      matches = false;\n      break;\n    }\n  }\n  if (matches) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not locate the skills paragraph block to reorder.\");\n}\n\nfor (let j = 0; j < newOrder.length; j++) {\n  if (oldOrder[j] !== newOrder[j]) {\n    paras.items[startIndex + j].insertText(newOrder[j], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The \"COMPETENCES TECHNIQUES\" skill list got reordered. The five skill\n# lines keep their exact text and paragraph formatting; only the order in\n# which they appear changes:\n#   old: Bases de donn\u00e9es, Autres, Visualisation, ML/AI, MLOps\n#   new: Visualisation, MLOps, Autres, ML/AI, Bases de donn\u00e9es\n# Since every one of those paragraphs shares identical paragraph\n# properties, the reorder can be applied by rewriting each paragraph's\n# text in place (cheapest, most robust transform -- no node moves needed).\n\n$d = $word.ActiveDocument\n\n$oldOrder = @(\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Autres : utilisables, 09\",\n  \"Visualisation : excel, tableau\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n)\n$newOrder = @(\n  \"Visualisation : excel, tableau\",\n  \"MLOps : vba, powerbi, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Autres : utilisables, 09\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n)\n\n$n = $oldOrder.Length\n$count = $d.Paragraphs.Count\n$startIndex = -1\n\nfor ($i = 1; $i -le ($count - $n + 1); $i++) {\n    $isMatch = $true\n    for ($j = 0; $j -lt $n; $j++) {\n        $ptext = $d.Paragraphs.Item($i + $j).Range.Text\n        $ptext = $ptext.TrimEnd([char]13, [char]7)\n        if ($ptext -ne $oldOrder[$j]) {\n            $isMatch = $false\n            break\n        }\n    }\n    if ($isMatch) {\n        $startIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1) {\n    throw \"Could not locate the skills paragraph block to reorder.\"\n}\n\nfor ($j = 0; $j -lt $n; $j++) {\n    if ($oldOrder[$j] -ne $newOrder[$j]) {\n        $d.Paragraphs.Item($startIndex + $j).Range.Text = $newOrder[$j]\n    }\n}\n"}
